$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 272 (shifts user.parent.id and everything below it down by one),
# matching where "user.locked" (row 271) leaves off alphabetically before "user.parent.id".
$ws.Rows.Item(272).Insert()
$ws.Rows.Item(272).RowHeight = 16.5

# Populate the new row 272 with the new search criterion: user.lockedUntrusted
$ws.Cells.Item(272, 1).Value = "USER_DB"
$ws.Cells.Item(272, 2).Value = "user.lockedUntrusted"
$ws.Cells.Item(272, 3).Value = "BOOLEAN"
$ws.Cells.Item(272, 12).Value = "user.lockedUntrusted"
$ws.Cells.Item(272, 13).Value = "EQ, NE"

# Reflect the scrolled/selected view state from the edit session.
$win = $excel.ActiveWindow
$win.ScrollRow = 262
$win.ScrollColumn = 1
$ws.Range("A272").Select()
